$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 45, shifting existing rows 45-65 down to 46-66.
$ws.Rows.Item(45).Insert()

# Populate the newly inserted row 45 with the new weekly record.
$ws.Cells.Item(45, 1).Value = 6
$ws.Cells.Item(45, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(45, 3).Value = "Metropolitana"
$ws.Cells.Item(45, 4).Value = 44609
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(45, 5).Value = 13
$ws.Cells.Item(45, 6).Value = 100114007
$ws.Cells.Item(45, 7).Value = "Jengibre"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 80
$ws.Cells.Item(45, 11).Value = 16000
$ws.Cells.Item(45, 12).Value = 17000
$ws.Cells.Item(45, 13).Value = 16375
$ws.Cells.Item(45, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(45, 15).Value = "Perú"
$ws.Cells.Item(45, 16).Value = 1260
$ws.Cells.Item(45, 17).Value = 13
$ws.Cells.Item(45, 18).Value = "Hortaliza"
